# Update the underlying worksheet data with new values (new performance unit).
# The 4 charts on the sheet reference these same ranges (Sheet1!$B$2:$G$8 and
# Sheet1!$B$14:$G$20), so updating the cell values will make Excel refresh the
# charts' cached numbers automatically on save.

function Set-RowValues($ws, $row, $values) {
    $n = $values.Length
    $arr = New-Object 'object[,]' 1, $n
    foreach ($j in 0..($n - 1)) {
        $arr[0, $j] = $values[$j]
    }
    $addr = "B" + $row + ":G" + $row
    $ws.Range($addr).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Table 1: rows 2-8 (columns B..G) ---
$table1 = @(
    @(13.796, 1182, 1716.689, 4193.725, 7213.091, 11652.414),
    @(20.967, 1269.41, 2516.29, 4942.29, 6114.948, 9467.282),
    @(20.258, 994.386, 1434.176, 2743.197, 4728.411, 12352.483),
    @(13.676, 1219.479, 2320.51, 3484.104, 7127.727, 13849.421),
    @(13.413, 731.642, 1515.003, 3419.862, 5521.32, 13419.69),
    @(20.686, 1118.107, 2288.893, 2516.611, 4736.65, 9700.235),
    @(12.488, 725.101, 1448.357, 3092.812, 6568.801, 11984.257)
)

for ($i = 0; $i -lt $table1.Length; $i++) {
    $rowNum = 2 + $i
    $rowVals = $table1[$i]
    Set-RowValues $ws $rowNum $rowVals
}

# --- Table 2: rows 14-20 (columns B..G) ---
$table2 = @(
    @(17.284, 1245.748, 1471.805, 2892.217, 4771.471, 10259.801),
    @(16.57, 723.47, 1887.367, 4359.349, 5217.945, 10002.001),
    @(18.933, 1042.739, 2069.397, 3064.457, 5094.694, 9275.524),
    @(18.504, 1265.333, 1940.981, 4940.194, 5162.239, 8909.114),
    @(16.594, 863.966, 2168.294, 3192.945, 5323.639, 9204.333),
    @(13.355, 1006.203, 1788.921, 4288.269, 5184.624, 9324.478),
    @(12.901, 897.935, 2385.285, 3854.598, 7212.185, 8906.049)
)

for ($i = 0; $i -lt $table2.Length; $i++) {
    $rowNum = 14 + $i
    $rowVals = $table2[$i]
    Set-RowValues $ws $rowNum $rowVals
}

# Leave the active selection on G14, matching the final state of the sheet.
$ws.Range("G14").Select()

$wb.Save()
